$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, and whether the value must be
# force-typed as text (leading apostrophe) because it would otherwise be
# auto-parsed by Excel as a number (losing significant trailing zeros).
$updates = @(
    ,@("D2", '27.883.04', $false)
    ,@("E2", '  -2.25%  ', $false)
    ,@("D3", '1.866.93', $false)
    ,@("E3", '  -2.37%  ', $false)
    ,@("E4", '  +0.24%  ', $false)
    ,@("D5", '311.98', $true)
    ,@("E5", '  -1.03%  ', $false)
    ,@("E6", '  +0.17%  ', $false)
    ,@("D7", '0.4957', $true)
    ,@("E7", '  -3.67%  ', $false)
    ,@("D8", '0.3795', $true)
    ,@("E8", '  -4.47%  ', $false)
    ,@("D9", '0.08884', $true)
    ,@("E9", '  -9.30%  ', $false)
    ,@("D10", '1.113', $true)
    ,@("E10", '  -3.29%  ', $false)
    ,@("D11", '41.56', $true)
    ,@("E11", '  -1.65%  ', $false)
    ,@("D12", '6.291', $true)
    ,@("E12", '  -3.78%  ', $false)
    ,@("D13", '20.57', $true)
    ,@("E13", '  -2.94%  ', $false)
    ,@("D14", '1.879.07', $false)
    ,@("E14", '  -1.30%  ', $false)
    ,@("D15", '7.197', $true)
    ,@("E15", '  -3.93%  ', $false)
    ,@("D16", '1.004', $true)
    ,@("E16", '  +0.29%  ', $false)
    ,@("D17", '0.00001095', $true)
    ,@("E17", '  -3.67%  ', $false)
    ,@("D18", '90.50', $true)
    ,@("E18", '  -4.35%  ', $false)
    ,@("D19", '0.06625', $true)
    ,@("E19", '  -0.52%  ', $false)
    ,@("D20", '17.82', $true)
    ,@("E20", '  -2.47%  ', $false)
    ,@("E21", '  +0.19%  ', $false)
    ,@("D22", '6.090', $true)
    ,@("E22", '  -3.71%  ', $false)
    ,@("D23", '27.931.33', $false)
    ,@("E23", '  -2.28%  ', $false)
    ,@("D24", '11.32', $true)
    ,@("E24", '  -1.58%  ', $false)
    ,@("D25", '2.282', $true)
    ,@("E25", '  -1.63%  ', $false)
    ,@("B26", 'LEO', $false)
    ,@("C26", 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo', $false)
    ,@("D26", '3.376', $true)
    ,@("E26", '  -0.05%  ', $false)
    ,@("B27", 'WrappedliquidstakedEther2.0', $false)
    ,@("C27", 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', $false)
    ,@("D27", '2.087.95', $false)
    ,@("E27", '  -1.84%  ', $false)
    ,@("B28", 'LidoDAOToken', $false)
    ,@("C28", 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', $false)
    ,@("D28", '2.497', $true)
    ,@("E28", '  -6.82%  ', $false)
    ,@("B29", 'Monero', $false)
    ,@("C29", 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', $false)
    ,@("D29", '157.68', $true)
    ,@("E29", '  +0.02%  ', $false)
    ,@("B30", 'EthereumClassic', $false)
    ,@("C30", 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', $false)
    ,@("D30", '20.69', $true)
    ,@("E30", '  -2.71%  ', $false)
    ,@("B31", 'BitcoinCash', $false)
    ,@("C31", 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', $false)
    ,@("D31", '126.49', $true)
    ,@("E31", '  -1.93%  ', $false)
    ,@("B32", 'Stellar', $false)
    ,@("C32", 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', $false)
    ,@("D32", '0.1052', $true)
    ,@("E32", '  -2.54%  ', $false)
    ,@("B33", 'ImmutableX', $false)
    ,@("C33", 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', $false)
    ,@("D33", '1.050', $true)
    ,@("E33", '  -5.85%  ', $false)
    ,@("B34", 'Filecoin', $false)
    ,@("C34", 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', $false)
    ,@("D34", '5.557', $true)
    ,@("E34", '  -3.53%  ', $false)
    ,@("B35", 'HuobiToken', $false)
    ,@("C35", 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', $false)
    ,@("D35", '3.591', $true)
    ,@("E35", '  -1.07%  ', $false)
    ,@("B36", 'FraxShare', $false)
    ,@("C36", 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', $false)
    ,@("D36", '9.264', $true)
    ,@("E36", '  -6.44%  ', $false)
    ,@("B37", 'Hedera', $false)
    ,@("C37", 'https://coinranking.com/coin/jad286TjB+hedera-hbar', $false)
    ,@("D37", '0.06510', $true)
    ,@("E37", '  -4.16%  ', $false)
    ,@("B38", 'VeChain', $false)
    ,@("C38", 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', $false)
    ,@("D38", '0.02392', $true)
    ,@("E38", '  -1.75%  ', $false)
    ,@("B39", 'Algorand', $false)
    ,@("C39", 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', $false)
    ,@("D39", '0.2173', $true)
    ,@("E39", '  -1.91%  ', $false)
    ,@("B40", 'TrustWalletToken', $false)
    ,@("C40", 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', $false)
    ,@("D40", '1.265', $true)
    ,@("E40", '  +6.39%  ', $false)
    ,@("B41", 'ARBITRUM', $false)
    ,@("C41", 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', $false)
    ,@("D41", '1.191', $true)
    ,@("E41", '  -6.49%  ', $false)
    ,@("B42", 'Aptos', $false)
    ,@("C42", 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', $false)
    ,@("D42", '11.62', $true)
    ,@("E42", '  -1.72%  ', $false)
    ,@("B43", 'TheSandbox', $false)
    ,@("C43", 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', $false)
    ,@("D43", '0.6332', $true)
    ,@("E43", '  -2.33%  ', $false)
    ,@("B44", 'InternetComputer(DFINITY)', $false)
    ,@("C44", 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', $false)
    ,@("D44", '4.866', $true)
    ,@("E44", '  -4.48%  ', $false)
    ,@("B45", 'Frax', $false)
    ,@("C45", 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', $false)
    ,@("D45", '1.002', $true)
    ,@("E45", '  +0.15%  ', $false)
    ,@("B46", 'EnergySwap', $false)
    ,@("C46", 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', $false)
    ,@("D46", '13.14', $true)
    ,@("E46", '  -3.52%  ', $false)
    ,@("B47", 'Decentraland', $false)
    ,@("C47", 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', $false)
    ,@("D47", '0.5962', $true)
    ,@("E47", '  -2.35%  ', $false)
    ,@("B48", 'WEMIXTOKEN', $false)
    ,@("C48", 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', $false)
    ,@("D48", '1.287', $true)
    ,@("E48", '  -0.46%  ', $false)
    ,@("B49", 'PancakeSwap', $false)
    ,@("C49", 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', $false)
    ,@("D49", '3.672', $true)
    ,@("E49", '  -2.92%  ', $false)
    ,@("B50", 'EOS', $false)
    ,@("C50", 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos', $false)
    ,@("D50", '1.208', $true)
    ,@("E50", '  +0.30%  ', $false)
    ,@("B51", 'NEARProtocol', $false)
    ,@("C51", 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', $false)
    ,@("D51", '1.958', $true)
    ,@("E51", '  -4.13%  ', $false)
)

foreach ($u in $updates) {
    $ref = $u[0]
    $value = $u[1]
    $forceText = $u[2]
    if ($forceText) {
        $ws.Range($ref).Value = "'" + $value
    } else {
        $ws.Range($ref).Value = $value
    }
}
